$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ------------------------------------------------------------------
$title = $d.Paragraphs.Item(1)
$titleRange = $title.Range
$titleRange.Collapse(0)
$titleRange.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaRange = $metaPara.Range
$metaRange.Collapse(1)
$metaRange.Text = "Meta description: Read a review of Cash Bunny, a fun cartoon farm-themed slot game with Wild and Bonus symbols. Try it for free and enjoy its Double Spin feature."

# Bold just the "Meta description" label (leave the rest of the
# sentence, starting at the colon, in regular weight).
$labelRange = $d.Range($metaPara.Range.Start, $metaPara.Range.Start + 16)
$labelRange.Bold = 1

# ------------------------------------------------------------------
# 2) Remove the duplicated "Play Cash Bunny for Free - Review" bold
#    paragraph that used to sit near the end of the document (right
#    before the italic meta-description paragraph). Skip the real
#    Heading1 title at the very top of the document - only the
#    plain-body duplicate near the bottom should be removed.
# ------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    if (($para.Range.Text -eq "Play Cash Bunny for Free - Review`r") -and ($para.Style.NameLocal -ne "Heading 1")) {
        $para.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph with the new
#    DALL-E prompt text, keeping its existing italic formatting.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bodyRange = $lastPara.Range.Duplicate()
$bodyRange.End = $bodyRange.End - 1
$bodyRange.Text = "Prompt: ""Create a feature image for the game Cash Bunny in a cartoon style with a happy Maya warrior wearing glasses"" For the feature image of Cash Bunny, DALLE could draw a cartoon image of a happy Maya warrior wearing glasses surrounded by the adorable animals on Old McDonald's farm. The image could have a bright color scheme to appeal to players and convey a cheerful atmosphere. The Maya warrior could be holding a carrot, which is the Bonus symbol in the game, while the animals could be shown happily frolicking around in the background. The image could also include the Cash Bunny game logo prominently displayed to make it instantly recognizable to players."

Write-Host "Done"
